# Restructure the "dissolved" sheet: move the CITATION/ROCK TYPE columns
# (B:C) to the end of the data block (R:S), shifting the chemistry columns
# (D:S) two places to the left (into B:Q). Applies to rows 1-18 (header +
# 17 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dissolved")

$lastRow = 18
$width = 18   # columns B..S inclusive

for ($r = 1; $r -le $lastRow; $r++) {
    $rangeAddr = "B" + $r + ":S" + $r
    $srcRange = $ws.Range($rangeAddr)
    $vals = $srcRange.Value2

    $rotated = New-Object 'object[,]' 1,$width

    # old D..S (source cols 3..18, 1-based within the 18-wide block) -> new B..Q (0..15)
    for ($i = 3; $i -le $width; $i++) {
        $rotated[0, $i - 3] = $vals[1, $i]
    }
    # old B..C (source cols 1..2) -> new R..S (16..17)
    $rotated[0, $width - 2] = $vals[1, 1]
    $rotated[0, $width - 1] = $vals[1, 2]

    $srcRange.Value = $rotated
}
